# "Last version before add more things"
# - Remove the PREGUNTAS sheet entirely (and its now-unused shared strings/styles).
# - CURSOS becomes/stays the single active sheet (activeTab -> 0, tabSelected -> true).
# - Set CURSOS!A4 to 12307.
# - Move the selection on CURSOS to B6.

$wb = $excel.ActiveWorkbook

# Drop the second sheet (PREGUNTAS); CURSOS remains the only, active sheet.
$wb.Worksheets.Item("PREGUNTAS").Delete()

$ws = $wb.Worksheets.Item("CURSOS")

# Fill in the new data row.
$ws.Range("A4").Value = 12307

# Update the selection/active cell shown when the sheet is opened.
$ws.Range("B6").Select()
